# Auto-generated update of Siren_Profits workbook market data
# Applies value changes + a few cell deletions across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3749.9167
$ws.Range("I32").Value = 2900
$ws.Range("K32").Value = 2900
$ws.Range("M32").Value = -2574
$ws.Range("H34").Value = 4400
$ws.Range("I34").Value = 4400
$ws.Range("K34").Value = 4400
$ws.Range("M34").Value = -4197
$ws.Range("H36").Value = 4400
$ws.Range("I36").Value = 4400
$ws.Range("K36").Value = 4400
$ws.Range("M36").Value = -3685
$ws.Range("H116").Value = 10188659
$ws.Range("J116").Value = 4183.857
$ws.Range("L116").Value = 4183.857
$ws.Range("N116").Value = -11067.857
$ws.Range("H135").Value = 6562.56
$ws.Range("I135").Value = 9325.467
$ws.Range("J135").Value = 2418.2
$ws.Range("K135").Value = 83929.20300000001
$ws.Range("L135").Value = 21763.8
$ws.Range("M135").Value = -81394.20300000001
$ws.Range("N135").Value = -26833.8
$ws.Range("H138").Value = 427463.7
$ws.Range("J138").Value = 4458
$ws.Range("L138").Value = 13374
$ws.Range("N138").Value = -23654
$ws.Range("H141").Value = 10749.333
$ws.Range("I141").Value = 10843
$ws.Range("K141").Value = 32529
$ws.Range("M141").Value = -27349

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3344.0386
$ws.Range("I32").Value = 3344.0386
$ws.Range("K32").Value = 3344.0386
$ws.Range("M32").Value = -3057.0386
$ws.Range("H45").Value = 190084.9
$ws.Range("I45").Value = 293279.84
$ws.Range("J45").Value = 9493.75
$ws.Range("K45").Value = 293279.84
$ws.Range("L45").Value = 9493.75
$ws.Range("M45").Value = -292902.84
$ws.Range("N45").Value = -10247.75
$ws.Range("H81").Value = 64999.5
$ws.Range("J81").Value = 64999.5
$ws.Range("L81").Value = 64999.5
$ws.Range("N81").Value = -66995.5
$ws.Range("H84").Value = 64999.5
$ws.Range("J84").Value = 64999.5
$ws.Range("L84").Value = 194998.5
$ws.Range("N84").Value = -204982.5
$ws.Range("H102").Value = 9962.5
$ws.Range("I102").Value = 13352.667
$ws.Range("J102").Value = 2334.625
$ws.Range("K102").Value = 13352.667
$ws.Range("L102").Value = 2334.625
$ws.Range("M102").Value = -11730.667
$ws.Range("N102").Value = -5578.625
$ws.Range("H122").Value = 1312668.6
$ws.Range("I122").Value = 9125.818
$ws.Range("K122").Value = 27377.454
$ws.Range("M122").Value = -24927.454

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3424.2727
$ws.Range("I20").Value = 2274.8572
$ws.Range("J20").Value = 5435.75
$ws.Range("K20").Value = 2274.8572
$ws.Range("L20").Value = 5435.75
$ws.Range("M20").Value = -2027.8572
$ws.Range("N20").Value = -5929.75
$ws.Range("H86").Value = 4830.2593
$ws.Range("I86").Value = 6135.9414
$ws.Range("J86").Value = 2610.6
$ws.Range("K86").Value = 6135.9414
$ws.Range("L86").Value = 2610.6
$ws.Range("M86").Value = -5012.9414
$ws.Range("N86").Value = -4856.6
$ws.Range("H89").Value = 4830.2593
$ws.Range("I89").Value = 6135.9414
$ws.Range("J89").Value = 2610.6
$ws.Range("K89").Value = 30679.707
$ws.Range("L89").Value = 13053
$ws.Range("M89").Value = -25063.707
$ws.Range("N89").Value = -24285

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8312.238
$ws.Range("I31").Value = 10261.357
$ws.Range("K31").Value = 10261.357
$ws.Range("M31").Value = -9966.357
$ws.Range("H34").Value = 8312.238
$ws.Range("I34").Value = 10261.357
$ws.Range("K34").Value = 10261.357
$ws.Range("M34").Value = -10059.357
$ws.Range("H87").Value = 53498.25
$ws.Range("J87").Value = 53498.25
$ws.Range("L87").Value = 53498.25
$ws.Range("N87").Value = -55870.25
$ws.Range("H90").Value = 53498.25
$ws.Range("J90").Value = 53498.25
$ws.Range("L90").Value = 160494.75
$ws.Range("N90").Value = -172350.75
$ws.Range("H99").Value = 150090.03
$ws.Range("I99").Value = 279846.06
$ws.Range("J99").Value = 4114.5
$ws.Range("K99").Value = 279846.06
$ws.Range("L99").Value = 4114.5
$ws.Range("M99").Value = -278348.06
$ws.Range("N99").Value = -7110.5
$ws.Range("H122").Value = 11073.75
$ws.Range("H126").Value = 150090.03
$ws.Range("I126").Value = 279846.06
$ws.Range("J126").Value = 4114.5
$ws.Range("K126").Value = 839538.1799999999
$ws.Range("L126").Value = 12343.5
$ws.Range("M126").Value = -837068.1799999999
$ws.Range("N126").Value = -17283.5
$ws.Range("H132").Value = 1546.1428
$ws.Range("J132").Value = 1554.6666
$ws.Range("L132").Value = 4663.9998
$ws.Range("N132").Value = -9723.9998
$ws.Range("H141").Value = 321930.47
$ws.Range("J141").Value = 356854.53
$ws.Range("L141").Value = 356854.53
$ws.Range("N141").Value = -367214.53

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 248.09091
$ws.Range("I8").Value = 248.09091
$ws.Range("K8").Value = 744.27273
$ws.Range("M8").Value = -605.27273
$ws.Range("H109").Value = 2087.5
$ws.Range("I109").Value = 1116.6666
$ws.Range("K109").Value = 3349.9998
$ws.Range("M109").Value = -2309.9998
$ws.Range("H131").Value = 1954.5333
$ws.Range("I131").Value = 1298.5
$ws.Range("J131").Value = 2018.5366
$ws.Range("K131").Value = 3895.5
$ws.Range("L131").Value = 6055.6098
$ws.Range("M131").Value = 1144.5
$ws.Range("N131").Value = -16135.6098
$ws.Range("H133").Value = 7966
$ws.Range("I133").Value = 7966
$ws.Range("K133").Value = 23898
$ws.Range("M133").Value = -18838

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8222.35
$ws.Range("I80").Value = 10158.23
$ws.Range("J80").Value = 4627.143
$ws.Range("K80").Value = 10158.23
$ws.Range("L80").Value = 4627.143
$ws.Range("M80").Value = -9160.23
$ws.Range("N80").Value = -6623.143
$ws.Range("H83").Value = 8222.35
$ws.Range("I83").Value = 10158.23
$ws.Range("J83").Value = 4627.143
$ws.Range("K83").Value = 50791.14999999999
$ws.Range("L83").Value = 23135.715
$ws.Range("M83").Value = -45799.14999999999
$ws.Range("N83").Value = -33119.715
$ws.Range("H102").Value = 10573.55
$ws.Range("I102").Value = 12728.267
$ws.Range("J102").Value = 4109.4
$ws.Range("K102").Value = 12728.267
$ws.Range("L102").Value = 4109.4
$ws.Range("M102").Value = -11106.267
$ws.Range("N102").Value = -7353.4
$ws.Range("H132").Value = 4076.3157
$ws.Range("I132").Value = 4254.355
$ws.Range("J132").Value = 3287.8572
$ws.Range("K132").Value = 12763.065
$ws.Range("L132").Value = 9863.5716
$ws.Range("M132").Value = -10233.065
$ws.Range("N132").Value = -14923.5716

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 134285710
$ws.Range("J2").Value = 18000000
$ws.Range("L2").Value = 18000000
$ws.Range("N2").Value = -18000224
$ws.Range("H7").Value = 24930.096
$ws.Range("I7").Value = 33395.285
$ws.Range("K7").Value = 33395.285
$ws.Range("M7").Value = -33283.285
$ws.Range("H40").Value = 17687.814
$ws.Range("I40").Value = 18279.9
$ws.Range("K40").Value = 18279.9
$ws.Range("M40").Value = -18143.9
$ws.Range("H93").Value = 5342.2104
$ws.Range("I93").Value = 6548.4287
$ws.Range("J93").Value = 1964.8
$ws.Range("K93").Value = 6548.4287
$ws.Range("L93").Value = 1964.8
$ws.Range("M93").Value = -5300.4287
$ws.Range("N93").Value = -4460.8
$ws.Range("H96").Value = 48000
$ws.Range("J96").Value = 48000
$ws.Range("L96").Value = 48000
$ws.Range("N96").Value = -53492
$ws.Range("H122").Value = 6569.8076
$ws.Range("I122").Value = 7293.1816
$ws.Range("J122").Value = 6039.3335
$ws.Range("K122").Value = 21879.5448
$ws.Range("L122").Value = 18118.0005
$ws.Range("M122").Value = -19429.5448
$ws.Range("N122").Value = -23018.0005
$ws.Range("H126").Value = 24930.096
$ws.Range("I126").Value = 33395.285
$ws.Range("K126").Value = 100185.855
$ws.Range("M126").Value = -97715.85500000001
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("H132").Value = 295404.56
$ws.Range("I132").Value = 440828.03
$ws.Range("J132").Value = 4557.5884
$ws.Range("K132").Value = 1322484.09
$ws.Range("L132").Value = 13672.7652
$ws.Range("M132").Value = -1319954.09
$ws.Range("N132").Value = -18732.7652
$ws.Range("H136").Value = 4851.769
$ws.Range("I136").Value = 2219.8572
$ws.Range("J136").Value = 5821.421
$ws.Range("K136").Value = 6659.571599999999
$ws.Range("L136").Value = 17464.263
$ws.Range("M136").Value = -4109.571599999999
$ws.Range("N136").Value = -22564.263
$ws.Range("N130").ClearContents()

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("H122").Value = 3649.349
$ws.Range("I122").Value = 2354.35
$ws.Range("J122").Value = 5901.522
$ws.Range("K122").Value = 7063.049999999999
$ws.Range("L122").Value = 17704.566
$ws.Range("M122").Value = -4613.049999999999
$ws.Range("N122").Value = -22604.566
$ws.Range("H126").Value = 18948.482
$ws.Range("I126").Value = 24041.1
$ws.Range("J126").Value = 4398.143
$ws.Range("K126").Value = 72123.29999999999
$ws.Range("L126").Value = 13194.429
$ws.Range("M126").Value = -69653.29999999999
$ws.Range("N126").Value = -18134.429
$ws.Range("H132").Value = 5652.7144
$ws.Range("I132").Value = 6158.1665
$ws.Range("K132").Value = 18474.4995
$ws.Range("M132").Value = -15944.4995
$ws.Range("N21").ClearContents()
$ws.Range("N35").ClearContents()

Write-Host "Done applying Siren_Profits updates"
